$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("9.81")

$ws.Range("A10").Value = "conditional(add(add(x, vel), add(add(x, x), conditional(x, x))), add(add(add(x, y), add(vel, x)), add(add(x, y), y)))"
$ws.Range("B10").Value = -419
